# Update res_bus/vm_pu.xlsx results for "case with 380 kV" (Case_1_228):
# bus 0's voltage setpoint (col B) drops from 1.05 pu to 1.02 pu, and all the
# per-bus voltage-magnitude results (cols C-F, I-N) for rows 2-25 are
# recomputed accordingly. Column H has no data and column G/A stay 1/index.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.006795699771644
$ws.Cells.Item(2, 4).Value = 1.030263256528299
$ws.Cells.Item(2, 5).Value = 1.009532559647175
$ws.Cells.Item(2, 6).Value = 1.00497242631925
$ws.Cells.Item(2, 9).Value = 1.02955528997503
$ws.Cells.Item(2, 10).Value = 1.012071080893318
$ws.Cells.Item(2, 11).Value = 1.033074774661431
$ws.Cells.Item(2, 12).Value = 1.012405490914704
$ws.Cells.Item(2, 13).Value = 1.007859237321635
$ws.Cells.Item(2, 14).Value = 1.013508337914362
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.008197251359901
$ws.Cells.Item(3, 4).Value = 1.030606612076188
$ws.Cells.Item(3, 5).Value = 1.010735825058346
$ws.Cells.Item(3, 6).Value = 1.007032313929777
$ws.Cells.Item(3, 9).Value = 1.029546965970127
$ws.Cells.Item(3, 10).Value = 1.013101395776095
$ws.Cells.Item(3, 11).Value = 1.033227698754684
$ws.Cells.Item(3, 12).Value = 1.01341151132103
$ws.Cells.Item(3, 13).Value = 1.009718434850895
$ws.Cells.Item(3, 14).Value = 1.014540115962452
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.009102639152691
$ws.Cells.Item(4, 4).Value = 1.030828240658591
$ws.Cells.Item(4, 5).Value = 1.011513358573906
$ws.Cells.Item(4, 6).Value = 1.008363069200797
$ws.Cells.Item(4, 9).Value = 1.029539753632458
$ws.Cells.Item(4, 10).Value = 1.013766239260105
$ws.Cells.Item(4, 11).Value = 1.0333253549493
$ws.Cells.Item(4, 12).Value = 1.014060866618908
$ws.Cells.Item(4, 13).Value = 1.010918994503388
$ws.Cells.Item(4, 14).Value = 1.015205903600468
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.009482911626778
$ws.Cells.Item(5, 4).Value = 1.030921280219724
$ws.Cells.Item(5, 5).Value = 1.011839986557137
$ws.Cells.Item(5, 6).Value = 1.008922028095641
$ws.Cells.Item(5, 9).Value = 1.029536283262584
$ws.Cells.Item(5, 10).Value = 1.014045305955706
$ws.Cells.Item(5, 11).Value = 1.033366097555123
$ws.Cells.Item(5, 12).Value = 1.014333476974628
$ws.Cells.Item(5, 13).Value = 1.011423136912126
$ws.Cells.Item(5, 14).Value = 1.015485366602786
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.009546740631775
$ws.Cells.Item(6, 4).Value = 1.030936894093018
$ws.Cells.Item(6, 5).Value = 1.011894814532448
$ws.Cells.Item(6, 6).Value = 1.009015851539715
$ws.Cells.Item(6, 9).Value = 1.029535674837201
$ws.Cells.Item(6, 10).Value = 1.014092137257115
$ws.Cells.Item(6, 11).Value = 1.033372920065691
$ws.Cells.Item(6, 12).Value = 1.014379227428651
$ws.Cells.Item(6, 13).Value = 1.011507751549646
$ws.Cells.Item(6, 14).Value = 1.015532264410014
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.009107721746463
$ws.Cells.Item(7, 4).Value = 1.030829484383945
$ws.Cells.Item(7, 5).Value = 1.011517723955896
$ws.Cells.Item(7, 6).Value = 1.008370539941099
$ws.Cells.Item(7, 9).Value = 1.029539708984861
$ws.Cells.Item(7, 10).Value = 1.013769969858014
$ws.Cells.Item(7, 11).Value = 1.033325900582478
$ws.Cells.Item(7, 12).Value = 1.014064510731977
$ws.Cells.Item(7, 13).Value = 1.010925733108616
$ws.Cells.Item(7, 14).Value = 1.015209639496253
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.007269678681445
$ws.Cells.Item(8, 4).Value = 1.030379406564799
$ws.Cells.Item(8, 5).Value = 1.009939432213309
$ws.Cells.Item(8, 6).Value = 1.005669027589001
$ws.Cells.Item(8, 9).Value = 1.02955285428455
$ws.Cells.Item(8, 10).Value = 1.012419665540274
$ws.Cells.Item(8, 11).Value = 1.033126723317437
$ws.Cells.Item(8, 12).Value = 1.012745816713894
$ws.Cells.Item(8, 13).Value = 1.008488082817075
$ws.Cells.Item(8, 14).Value = 1.013857417591499
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.004018848159921
$ws.Cells.Item(9, 4).Value = 1.029582245551909
$ws.Cells.Item(9, 5).Value = 1.007149868128996
$ws.Cells.Item(9, 6).Value = 1.000891439981142
$ws.Cells.Item(9, 9).Value = 1.029562081003761
$ws.Cells.Item(9, 10).Value = 1.010025865257892
$ws.Cells.Item(9, 11).Value = 1.032765901988229
$ws.Cells.Item(9, 12).Value = 1.010409523095615
$ws.Cells.Item(9, 13).Value = 1.00417296713567
$ws.Cells.Item(9, 14).Value = 1.011460217838144
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.001843021329063
$ws.Cells.Item(10, 4).Value = 1.029048220613477
$ws.Cells.Item(10, 5).Value = 1.005284085024044
$ws.Cells.Item(10, 6).Value = 0.9976936178203025
$ws.Cells.Item(10, 9).Value = 1.029558927936093
$ws.Cells.Item(10, 10).Value = 1.008419891959954
$ws.Cells.Item(10, 11).Value = 1.032518842347417
$ws.Cells.Item(10, 12).Value = 1.008843141418668
$ws.Cells.Item(10, 13).Value = 1.001281905842057
$ws.Cells.Item(10, 14).Value = 1.009851963873915
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.000898696951886
$ws.Cells.Item(11, 4).Value = 1.028816400681362
$ws.Cells.Item(11, 5).Value = 1.004474645911796
$ws.Cells.Item(11, 6).Value = 0.9963056270408939
$ws.Cells.Item(11, 9).Value = 1.029555369739793
$ws.Cells.Item(11, 10).Value = 1.00772199414136
$ws.Cells.Item(11, 11).Value = 1.032410340936442
$ws.Cells.Item(11, 12).Value = 1.008162693502325
$ws.Cells.Item(11, 13).Value = 1.000026403695264
$ws.Cells.Item(11, 14).Value = 1.009153074960369
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.000547594827383
$ws.Cells.Item(12, 4).Value = 1.028730207336787
$ws.Cells.Item(12, 5).Value = 1.004173744240184
$ws.Cells.Item(12, 6).Value = 0.9957895443653946
$ws.Cells.Item(12, 9).Value = 1.02955371971951
$ws.Cells.Item(12, 10).Value = 1.007462379703838
$ws.Cells.Item(12, 11).Value = 1.032369811708957
$ws.Cells.Item(12, 12).Value = 1.007909607544663
$ws.Cells.Item(12, 13).Value = 0.9995594842477433
$ws.Cells.Item(12, 14).Value = 1.008893091840568
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.000622922897103
$ws.Cells.Item(13, 4).Value = 1.028748699911505
$ws.Cells.Item(13, 5).Value = 1.004238299660311
$ws.Cells.Item(13, 6).Value = 0.9959002698105656
$ws.Cells.Item(13, 9).Value = 1.029554088495152
$ws.Cells.Item(13, 10).Value = 1.007518085399739
$ws.Cells.Item(13, 11).Value = 1.032378515607504
$ws.Cells.Item(13, 12).Value = 1.00796391072889
$ws.Cells.Item(13, 13).Value = 0.9996596661921244
$ws.Cells.Item(13, 14).Value = 1.008948876644946
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.000869681699599
$ws.Cells.Item(14, 4).Value = 1.028809277633772
$ws.Cells.Item(14, 5).Value = 1.004449778228142
$ws.Cells.Item(14, 6).Value = 0.9962629782525639
$ws.Cells.Item(14, 9).Value = 1.029555240040407
$ws.Cells.Item(14, 10).Value = 1.007700542215689
$ws.Cells.Item(14, 11).Value = 1.032406995397466
$ws.Cells.Item(14, 12).Value = 1.008141780263487
$ws.Cells.Item(14, 13).Value = 0.9999878197400796
$ws.Cells.Item(14, 14).Value = 1.009131592570504
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.001021672823927
$ws.Cells.Item(15, 4).Value = 1.028846590368057
$ws.Cells.Item(15, 5).Value = 1.004580045112463
$ws.Cells.Item(15, 6).Value = 0.9964863851390772
$ws.Cells.Item(15, 9).Value = 1.029555906069241
$ws.Cells.Item(15, 10).Value = 1.007812908730181
$ws.Cells.Item(15, 11).Value = 1.032424512714591
$ws.Cells.Item(15, 12).Value = 1.00825132662824
$ws.Cells.Item(15, 13).Value = 1.000189929797714
$ws.Cells.Item(15, 14).Value = 1.009244118658335
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.001905646230184
$ws.Cells.Item(16, 4).Value = 1.029063593663635
$ws.Cells.Item(16, 5).Value = 1.005337771633596
$ws.Cells.Item(16, 6).Value = 0.9977856624107103
$ws.Cells.Item(16, 9).Value = 1.029559117975007
$ws.Cells.Item(16, 10).Value = 1.008466155736214
$ws.Cells.Item(16, 11).Value = 1.032526011298701
$ws.Cells.Item(16, 12).Value = 1.008888253624534
$ws.Cells.Item(16, 13).Value = 1.00136515060923
$ws.Cells.Item(16, 14).Value = 1.009898293350044
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.002459548863428
$ws.Cells.Item(17, 4).Value = 1.029199559676678
$ws.Cells.Item(17, 5).Value = 1.005812654693122
$ws.Cells.Item(17, 6).Value = 0.99859976137455
$ws.Cells.Item(17, 9).Value = 1.029560546404202
$ws.Cells.Item(17, 10).Value = 1.00887524489243
$ws.Cells.Item(17, 11).Value = 1.032589272239751
$ws.Cells.Item(17, 12).Value = 1.00928718825863
$ws.Cells.Item(17, 13).Value = 1.002101342565271
$ws.Cells.Item(17, 14).Value = 1.010307963459785
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.002782421593638
$ws.Cells.Item(18, 4).Value = 1.029278809869972
$ws.Cells.Item(18, 5).Value = 1.006089497654567
$ws.Cells.Item(18, 6).Value = 0.9990742934944372
$ws.Cells.Item(18, 9).Value = 1.029561167926103
$ws.Cells.Item(18, 10).Value = 1.009113618959025
$ws.Cells.Item(18, 11).Value = 1.032626024211051
$ws.Cells.Item(18, 12).Value = 1.009519668882094
$ws.Cells.Item(18, 13).Value = 1.002530399804979
$ws.Cells.Item(18, 14).Value = 1.010546676044896
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.002892477678187
$ws.Cells.Item(19, 4).Value = 1.029305822478735
$ws.Cells.Item(19, 5).Value = 1.00618386904272
$ws.Cells.Item(19, 6).Value = 0.9992360434449757
$ws.Cells.Item(19, 9).Value = 1.029561343914426
$ws.Cells.Item(19, 10).Value = 1.009194857778418
$ws.Cells.Item(19, 11).Value = 1.032638530682611
$ws.Cells.Item(19, 12).Value = 1.009598903173331
$ws.Cells.Item(19, 13).Value = 1.002676638324288
$ws.Cells.Item(19, 14).Value = 1.010628030232731
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.002400142081356
$ws.Cells.Item(20, 4).Value = 1.029184977633816
$ws.Cells.Item(20, 5).Value = 1.005761719632777
$ws.Cells.Item(20, 6).Value = 0.9985124493247259
$ws.Cells.Item(20, 9).Value = 1.029560415031907
$ws.Cells.Item(20, 10).Value = 1.008831378456744
$ws.Cells.Item(20, 11).Value = 1.032582500137902
$ws.Cells.Item(20, 12).Value = 1.009244408248378
$ws.Cells.Item(20, 13).Value = 1.002022392601343
$ws.Cells.Item(20, 14).Value = 1.010264034728729
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.000797026801523
$ws.Cells.Item(21, 4).Value = 1.028791441332969
$ws.Cells.Item(21, 5).Value = 1.004387509750594
$ws.Cells.Item(21, 6).Value = 0.9961561842308329
$ws.Cells.Item(21, 9).Value = 1.029554909994161
$ws.Cells.Item(21, 10).Value = 1.007646823902336
$ws.Cells.Item(21, 11).Value = 1.032398615060048
$ws.Cells.Item(21, 12).Value = 1.00808941149181
$ws.Cells.Item(21, 13).Value = 0.9998912026063483
$ws.Cells.Item(21, 14).Value = 1.009077797970983
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 0.9997871235086685
$ws.Cells.Item(22, 4).Value = 1.028543518342299
$ws.Cells.Item(22, 5).Value = 1.003522097095154
$ws.Cells.Item(22, 6).Value = 0.9946716801406038
$ws.Cells.Item(22, 9).Value = 1.029549549582982
$ws.Cells.Item(22, 10).Value = 1.006899820848372
$ws.Cells.Item(22, 11).Value = 1.03228168691861
$ws.Cells.Item(22, 12).Value = 1.007361263806506
$ws.Cells.Item(22, 13).Value = 0.9985479300973373
$ws.Cells.Item(22, 14).Value = 1.008329734087001
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.000322681920053
$ws.Cells.Item(23, 4).Value = 1.028674992704237
$ws.Cells.Item(23, 5).Value = 1.003981003448065
$ws.Cells.Item(23, 6).Value = 0.9954589386253812
$ws.Cells.Item(23, 9).Value = 1.029552570889732
$ws.Cells.Item(23, 10).Value = 1.007296035183665
$ws.Cells.Item(23, 11).Value = 1.032343796520178
$ws.Cells.Item(23, 12).Value = 1.007747456589261
$ws.Cells.Item(23, 13).Value = 0.9992603450241445
$ws.Cells.Item(23, 14).Value = 1.008726511092096
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.002426986095816
$ws.Cells.Item(24, 4).Value = 1.029191566806252
$ws.Cells.Item(24, 5).Value = 1.005784735452404
$ws.Cells.Item(24, 6).Value = 0.9985519028639323
$ws.Cells.Item(24, 9).Value = 1.029560475047558
$ws.Cells.Item(24, 10).Value = 1.008851200554655
$ws.Cells.Item(24, 11).Value = 1.032585560613934
$ws.Cells.Item(24, 12).Value = 1.009263739346227
$ws.Cells.Item(24, 13).Value = 1.00205806777604
$ws.Cells.Item(24, 14).Value = 1.010283884976293
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.004860741150578
$ws.Cells.Item(25, 4).Value = 1.029788797171083
$ws.Cells.Item(25, 5).Value = 1.007872076959016
$ws.Cells.Item(25, 6).Value = 1.002128720824313
$ws.Cells.Item(25, 9).Value = 1.029561340424493
$ws.Cells.Item(25, 10).Value = 1.010646469523226
$ws.Cells.Item(25, 11).Value = 1.032860338570589
$ws.Cells.Item(25, 12).Value = 1.012081703432719
$ws.Cells.Item(25, 13).Value = 1.005290970498297
$ws.Cells.Item(25, 14).Value = 1.012081703432719
